$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, matching the original inline-string
# representation (values like "160.80" or "28.332.66" must not be reinterpreted
# as numbers by Excel, which would silently drop trailing zeros / merge the dots).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.303.45'
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("D3").Value = '1.811.11'
$ws.Range("E3").Value = '  -0.75%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").Value = '312.41'
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("D7").Value = '0.5161'
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = '0.3972'
$ws.Range("E8").Value = '  +3.07%  '
$ws.Range("D9").Value = '0.07828'
$ws.Range("E9").Value = '  -5.40%  '
$ws.Range("D10").Value = '1.112'
$ws.Range("E10").Value = '  -1.04%  '
$ws.Range("D11").Value = '41.05'
$ws.Range("E11").Value = '  -1.93%  '
$ws.Range("D12").Value = '6.351'
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D14").Value = '20.43'
$ws.Range("E14").Value = '  -3.03%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '7.310'
$ws.Range("E15").Value = '  -2.14%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.807.27'
$ws.Range("E16").Value = '  -0.57%  '
$ws.Range("D17").Value = '92.62'
$ws.Range("E17").Value = '  -1.52%  '
$ws.Range("E18").Value = '  -3.43%  '
$ws.Range("D19").Value = '0.06581'
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = '17.29'
$ws.Range("D22").Value = '6.007'
$ws.Range("E22").Value = '  -0.61%  '
$ws.Range("D23").Value = '28.332.66'
$ws.Range("E23").Value = '  -0.61%  '
$ws.Range("D24").Value = '11.13'
$ws.Range("E24").Value = '  -3.66%  '
$ws.Range("D25").Value = '2.228'
$ws.Range("E25").Value = '  -0.74%  '
$ws.Range("D26").Value = '160.80'
$ws.Range("E26").Value = '  +0.92%  '
$ws.Range("D27").Value = '2.424'
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("E28").Value = '  -2.63%  '
$ws.Range("D29").Value = '2.019.21'
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("D30").Value = '127.47'
$ws.Range("E30").Value = '  +1.44%  '
$ws.Range("E31").Value = '  -0.64%  '
$ws.Range("D32").Value = '1.062'
$ws.Range("D33").Value = '3.659'
$ws.Range("E33").Value = '  -0.91%  '
$ws.Range("E34").Value = '  -2.76%  '
$ws.Range("D35").Value = '0.07182'
$ws.Range("E35").Value = '  -4.55%  '
$ws.Range("D36").Value = '9.123'
$ws.Range("E36").Value = '  +4.01%  '
$ws.Range("D37").Value = '0.02360'
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("D38").Value = '0.2189'
$ws.Range("E38").Value = '  -1.70%  '
$ws.Range("D39").Value = '5.049'
$ws.Range("E39").Value = '  -3.43%  '
$ws.Range("E40").Value = '  -5.55%  '
$ws.Range("D41").Value = '0.6182'
$ws.Range("E41").Value = '  -3.22%  '
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("D44").Value = '13.29'
$ws.Range("E44").Value = '  -2.72%  '
$ws.Range("D45").Value = '0.6002'
$ws.Range("E45").Value = '  -3.26%  '
$ws.Range("D46").Value = '1.306'
$ws.Range("E46").Value = '  -6.41%  '
$ws.Range("D47").Value = '3.741'
$ws.Range("E47").Value = '  -1.65%  '
$ws.Range("D48").Value = '125.20'
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("D49").Value = '1.218'
$ws.Range("E49").Value = '  +1.15%  '
$ws.Range("E50").Value = '  -4.67%  '
$ws.Range("D51").Value = '0.06825'
$ws.Range("E51").Value = '  -1.77%  '
